# Implements "sample monitor UI" language-table additions.
# Appends rows 224-232 to Sheet1 (columns: A=index, B=name, C=en, D=zh_tw,
# E=de, F=zh_cn) for the new Sample Monitor / Test Setup / batch-config /
# move(Home,Next,Last) / Confirm UI strings, then leaves the active
# selection on the last-entered row (E232:F232), matching how the rows
# were typed into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- row 224 : run_stauts_monitor_title -------------------------------
$ws.Cells.Item(224, 1).Value = 221
$ws.Cells.Item(224, 2).Value = "run_stauts_monitor_title"
$ws.Cells.Item(224, 3).Value = "Sample Monitor"
$ws.Cells.Item(224, 4).Value = "樣品即時監測"
$ws.Cells.Item(224, 5).Value = "Probe-Monitor"
$ws.Cells.Item(224, 6).Value = "样品即时监测"

# --- row 225 : run_open_sampleSetupList_btn ----------------------------
$ws.Cells.Item(225, 1).Value = 222
$ws.Cells.Item(225, 2).Value = "run_open_sampleSetupList_btn"
$ws.Cells.Item(225, 3).Value = "Test Setup"
$ws.Cells.Item(225, 4).Value = "測試設定"
$ws.Cells.Item(225, 5).Value = "Versuchsaufbau"
$ws.Cells.Item(225, 6).Value = "测试设定"

# --- row 226 : modal_batch_setup_dialog_title --------------------------
$ws.Cells.Item(226, 1).Value = 223
$ws.Cells.Item(226, 2).Value = "modal_batch_setup_dialog_title"
$ws.Cells.Item(226, 3).Value = "Test Setup"
$ws.Cells.Item(226, 4).Value = "測試設定"
$ws.Cells.Item(226, 5).Value = "Versuchsaufbau"
$ws.Cells.Item(226, 6).Value = "测试设定"

# --- row 227 : sampleBatchConfigAdd_btn --------------------------------
$ws.Cells.Item(227, 1).Value = 224
$ws.Cells.Item(227, 2).Value = "sampleBatchConfigAdd_btn"
$ws.Cells.Item(227, 3).Value = "Add"
$ws.Cells.Item(227, 4).Value = "新增"
$ws.Cells.Item(227, 5).Value = "Hinzufügen"
$ws.Cells.Item(227, 6).Value = "新增"

# --- row 228 : sampleBatchConfigClearAll_btn ---------------------------
$ws.Cells.Item(228, 1).Value = 225
$ws.Cells.Item(228, 2).Value = "sampleBatchConfigClearAll_btn"
$ws.Cells.Item(228, 3).Value = "Clear All Samples"
$ws.Cells.Item(228, 4).Value = "清除所有樣品"
$ws.Cells.Item(228, 5).Value = "Alle Proben"
$ws.Cells.Item(228, 6).Value = "清除所有样品"

# --- row 229 : moveLastBtn (to last position) ---------------------------
$ws.Cells.Item(229, 1).Value = 226
$ws.Cells.Item(229, 2).Value = "moveLastBtn"
$ws.Cells.Item(229, 3).Value = "To Last Position"
$ws.Cells.Item(229, 4).Value = "移動至上個位置"
$ws.Cells.Item(229, 5).Value = "Letzte Position"
$ws.Cells.Item(229, 6).Value = "移动至上个位置"

# --- row 230 : moveHomeBtn (to home) -------------------------------------
$ws.Cells.Item(230, 1).Value = 227
$ws.Cells.Item(230, 2).Value = "moveHomeBtn"
$ws.Cells.Item(230, 3).Value = "To Home"
$ws.Cells.Item(230, 4).Value = "回原點"
$ws.Cells.Item(230, 5).Value = "Nach Hause"
$ws.Cells.Item(230, 6).Value = "回原点"

# --- row 231 : moveNextBtn (to next position) ----------------------------
$ws.Cells.Item(231, 1).Value = 228
$ws.Cells.Item(231, 2).Value = "moveNextBtn"
$ws.Cells.Item(231, 3).Value = "To Next Position"
$ws.Cells.Item(231, 4).Value = "移動至下個位置"
$ws.Cells.Item(231, 5).Value = "Zum nächsten Position"
$ws.Cells.Item(231, 6).Value = "移动至下个位置"

# --- row 232 : sampleSetupConfirmBtn (no index in column A) -------------
$ws.Cells.Item(232, 2).Value = "sampleSetupConfirmBtn"
$ws.Cells.Item(232, 3).Value = "Confirm"
$ws.Cells.Item(232, 4).Value = "確定"
$ws.Cells.Item(232, 5).Value = "Bestätigen"
$ws.Cells.Item(232, 6).Value = "确定"

# --- Match formatting of the rest of the index/name table ---------------
# Column A uses the bordered/bold "index" style used throughout the sheet
# (rows 2-223); copy it down onto the new index cells (224-231) without
# disturbing the values just written.
$ws.Cells.Item(223, 1).Copy()
$ws.Range($ws.Cells.Item(224, 1), $ws.Cells.Item(231, 1)).PasteSpecial($xlPasteFormats)

# D224 picked up the odd one-off font style also seen on D4/D67/D73/... ;
# replicate it so the new row matches.
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(224, 4).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Leave the selection where the author's last edit landed.
$ws.Range("E232:F232").Select() | Out-Null
